$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell with shared string "hi"
$ws.Range("D4").Value = "hi"

# Column A: simple fill of sequential integers (not formulas per diff, except A7 onward)
$ws.Range("A5").Value = 28
$ws.Range("A6").Value = 29
$ws.Range("A7").Formula = "=A6+1"
$ws.Range("A8:A17").Formula = "=A7+1"

# Column B formulas
$ws.Range("B5").Formula = "=B4-3"
$ws.Range("B6:B11").Formula = "=B5-3"
$ws.Range("B12").Formula = "=B11+2"
$ws.Range("B13:B17").Formula = "=B12+2"

# Update selection to D4 as in diff
$ws.Range("D4").Select()
